$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.088846
$ws.Range("H2").Value = 30.266538
$ws.Range("I2").Value = 0.1151445838515654
$ws.Range("J2").Value = 0.1151445838515654
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 1696.034159052591
$ws.Range("R2").Value = 15264.30743147332
$ws.Range("S2").Value = 0.03436124242383793
$ws.Range("T2").Value = 0.03436124242383794
$ws.Range("G3").Value = 10.088846
$ws.Range("H3").Value = 30.266538
$ws.Range("I3").Value = 0.1151445838515654
$ws.Range("J3").Value = 0.1151445838515654
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 1644.544808680707
$ws.Range("R3").Value = 14800.90327812637
$ws.Range("S3").Value = 0.03331808062138785
$ws.Range("T3").Value = 0.03331808062138785
$ws.Range("G4").Value = 10.088846
$ws.Range("H4").Value = 30.266538
$ws.Range("I4").Value = 0.1151445838515654
$ws.Range("J4").Value = 0.1151445838515654
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 1674.68316116638
$ws.Range("R4").Value = 15072.14845049742
$ws.Range("S4").Value = 0.03392867636351238
$ws.Range("T4").Value = 0.03392867636351238
$ws.Range("G5").Value = 10.088846
$ws.Range("H5").Value = 30.266538
$ws.Range("I5").Value = 0.1151445838515654
$ws.Range("J5").Value = 0.1151445838515654
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 668.1513237719128
$ws.Range("R5").Value = 6013.361913947214
$ws.Range("S5").Value = 0.01353658444282727
$ws.Range("T5").Value = 0.01353658444282727
$ws.Range("I6").Value = 0.4327250566572728
$ws.Range("J6").Value = 0.4327250566572729
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 6373.868861385655
$ws.Range("R6").Value = 57364.8197524709
$ws.Range("S6").Value = 0.129133043668275
$ws.Range("T6").Value = 0.1291330436682751
$ws.Range("I7").Value = 0.4327250566572728
$ws.Range("J7").Value = 0.4327250566572729
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.1252127355220419
$ws.Range("T7").Value = 0.125212735522042
$ws.Range("I8").Value = 0.4327250566572728
$ws.Range("J8").Value = 0.4327250566572729
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 6293.629639869947
$ws.Range("R8").Value = 56642.66675882952
$ws.Range("S8").Value = 0.1275074164203301
$ws.Range("T8").Value = 0.1275074164203301
$ws.Range("I9").Value = 0.4327250566572728
$ws.Range("J9").Value = 0.4327250566572729
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 2510.980627691087
$ws.Range("R9").Value = 22598.82564921979
$ws.Range("S9").Value = 0.05087186104662578
$ws.Range("T9").Value = 0.05087186104662578
$ws.Range("G10").Value = 15.69885766666667
$ws.Range("H10").Value = 47.096573
$ws.Range("I10").Value = 0.1791719719949428
$ws.Range("J10").Value = 0.1791719719949428
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 2639.13225167391
$ws.Range("R10").Value = 23752.19026506519
$ws.Range("S10").Value = 0.05346818199640078
$ws.Range("T10").Value = 0.05346818199640078
$ws.Range("G11").Value = 15.69885766666667
$ws.Range("H11").Value = 47.096573
$ws.Range("I11").Value = 0.1791719719949428
$ws.Range("J11").Value = 0.1791719719949428
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 2559.011692510123
$ws.Range("R11").Value = 23031.10523259111
$ws.Range("S11").Value = 0.05184495881904558
$ws.Range("T11").Value = 0.05184495881904558
$ws.Range("G12").Value = 15.69885766666667
$ws.Range("H12").Value = 47.096573
$ws.Range("I12").Value = 0.1791719719949428
$ws.Range("J12").Value = 0.1791719719949428
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 2605.908801057563
$ws.Range("R12").Value = 23453.17920951807
$ws.Range("S12").Value = 0.05279508291128426
$ws.Range("T12").Value = 0.05279508291128426
$ws.Range("G13").Value = 15.69885766666667
$ws.Range("H13").Value = 47.096573
$ws.Range("I13").Value = 0.1791719719949428
$ws.Range("J13").Value = 0.1791719719949428
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 1039.684076027147
$ws.Range("R13").Value = 9357.156684244319
$ws.Range("S13").Value = 0.0210637482682122
$ws.Range("T13").Value = 0.0210637482682122
$ws.Range("G14").Value = 23.91632366666667
$ws.Range("H14").Value = 71.748971
$ws.Range("I14").Value = 0.2729583874962189
$ws.Range("J14").Value = 0.2729583874962189
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 4020.569042051447
$ws.Range("R14").Value = 36185.12137846302
$ws.Range("S14").Value = 0.0814557577147382
$ws.Range("T14").Value = 0.08145575771473822
$ws.Range("G15").Value = 23.91632366666667
$ws.Range("H15").Value = 71.748971
$ws.Range("I15").Value = 0.2729583874962189
$ws.Range("J15").Value = 0.2729583874962189
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 3898.509891888944
$ws.Range("R15").Value = 35086.5890270005
$ws.Range("S15").Value = 0.07898286881306407
$ws.Range("T15").Value = 0.07898286881306407
$ws.Range("G16").Value = 23.91632366666667
$ws.Range("H16").Value = 71.748971
$ws.Range("I16").Value = 0.2729583874962189
$ws.Range("J16").Value = 0.2729583874962189
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 3969.954990052543
$ws.Range("R16").Value = 35729.59491047289
$ws.Range("S16").Value = 0.08043032924591624
$ws.Range("T16").Value = 0.08043032924591624
$ws.Range("G17").Value = 23.91632366666667
$ws.Range("H17").Value = 71.748971
$ws.Range("I17").Value = 0.2729583874962189
$ws.Range("J17").Value = 0.2729583874962189
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 1583.900013702346
$ws.Range("R17").Value = 14255.10012332111
$ws.Range("S17").Value = 0.03208943172250042
$ws.Range("T17").Value = 0.03208943172250042
